$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Drop the stray <w:lastRenderedPageBreak/> that sits in front of the run
# "time. Essentially, the price per lock ...". A self-replace (via Find)
# forces Word to rebuild the run without the stale render marker.
$d.Content.Find.Execute("time. Essentially, the price per lock", $true, $false, $false, $false, $false, $true, 1, $false, "time. Essentially, the price per lock", 2)

# --- Change 2 -------------------------------------------------------------
# The _GoBack bookmark moves from the end of answer "3." to the end of
# answer "5.". Remove it from its old spot first.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# Answer "5." currently just reads "5.  " (two trailing spaces). Replace it
# with the full answer text, keeping both trailing spaces - the bookmark
# will later be dropped in between them.
$d.Content.Find.Execute("5.  ", $true, $false, $false, $false, $false, $true, 1, $false, "5. This cannot be done with a user-mode implementation of pthread_cond_wait. This needs to be an atomic operation, and the only way to ensure that it’s an atomic operation is by using a kernel system call. In user mode, nothing is guaranteed, so the only safe way to do this is through kernel mode.  ", 2)

# Re-create _GoBack right before the final trailing space of answer "5.".
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p5.Range
$pos = $r.End - 2
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
